# Livres.xlsx — fix the swapped <categorie>/<image> tag formulas in the
# I column on sheets Feuil2, Feuil3, Feuil4 and Feuil6 (Feuil1 already has
# the correct order), bump the book-id seed on Feuil2 from 5 to 6, and
# update the active-sheet / selection state to match the author's last
# interaction (reading the "livres" table for display on Feuil1, after
# having walked through Feuil2..Feuil6).

$wb = $excel.ActiveWorkbook

# Each "book" record occupies 9 rows; within a record the two rows at
# offset +5/+6 from the record's first row hold the <categorie>/<image>
# formulas (built from the same A-column TRIM()) — on Feuil2..Feuil6 they
# were written in the wrong order.
$pairs = @(6, 15, 24, 33, 42)

$sheetNames = @("Feuil2", "Feuil3", "Feuil4", "Feuil6")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    foreach ($r1 in $pairs) {
        $r2 = $r1 + 1
        $ws.Range("I$r1").Formula = "=CONCATENATE(`"<image>`",TRIM(A$r1),`"</image>`")"
        $ws.Range("I$r2").Formula = "=CONCATENATE(`"<categorie>`",TRIM(A$r2),`"</categorie>`")"
    }
}

# Feuil2's book-id seed moves from 5 to 6 (A10/A19/A28/A37 and the H-column
# "<livre id=...>" strings are formulas and recompute automatically).
$ws2 = $wb.Worksheets.Item("Feuil2")
$ws2.Range("A1").Value = 6

# Walk the view state the same way the author did: select/scroll through
# Feuil2, Feuil3, Feuil4 and Feuil6, ending back on Feuil1 as the active
# sheet with I6:I7 selected.
$ws2.Range("A2").Select()

$ws3 = $wb.Worksheets.Item("Feuil3")
$ws3.Range("A3").Select()

$ws4 = $wb.Worksheets.Item("Feuil4")
$ws4.Range("A3").Select()

$ws6 = $wb.Worksheets.Item("Feuil6")
$ws6.Range("A3").Select()

$ws1 = $wb.Worksheets.Item("Feuil1")
$ws1.Activate()
$ws1.Range("I6:I7").Select()
